$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$layouts = $m.CustomLayouts
$lyt = $layouts.Item(1)
$sh = $lyt.Shapes.Item(5)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$w = $tr.Words(1,1)
$w.Text = "5/8/2024"
